# Apply edit: add columns I (I0) and J (IF) with header + values for rows 2-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), matching style of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Style = $ws.Range("H1").Style

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Style = $ws.Range("H1").Style

# Data values for I2:J28
$values = @(
    @(2, 6, 6),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 7, 7),
    @(6, 9, 9),
    @(7, 7, 7),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 7, 8),
    @(11, 1, 2),
    @(12, 7, 8),
    @(13, 7, 7),
    @(14, 6, 7),
    @(15, 6, 6),
    @(16, 5, 6),
    @(17, 8, 8),
    @(18, 6, 6),
    @(19, 6, 6),
    @(20, 7, 7),
    @(21, 9, 9),
    @(22, 8, 9),
    @(23, 9, 9),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 6, 6),
    @(28, 8, 8)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
